$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B should mirror column A for every used row (translation patch v2:
# duplicate the existing column A strings into the new column B).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Range("A$r").Copy($ws.Range("B$r"))
}
